# ----------------------------------------------------------------------------
# Add default cut-in / cut-out wind-speed parameter rows (rows 46-55) to the
# "Operation" category, for all turbine sizes, replacing values that used to
# be computed by the (now deleted) cut_speeds.py helper script.
# ----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-create the two new number-format styles in the same order Excel would
# (numFmt "0.0" before built-in integer numFmt "0"), so cellXfs indices land
# the same way they did in the authored workbook (idx 7 -> "0.0", idx 8 -> "0").
$ws.Range("K51").Value = 23.669871794871799
$ws.Range("K51").NumberFormat = "0.0"
$ws.Range("K51").HorizontalAlignment = -4108
$ws.Range("K51").VerticalAlignment = -4108

$ws.Range("C46").Value = 100
$ws.Range("C46").NumberFormat = "0"
$ws.Range("C46").HorizontalAlignment = -4108
$ws.Range("C46").VerticalAlignment = -4108

# --- Row 46 ---
$ws.Range("A46").Value = 'Operation'
$ws.Range("B46").Value = 'all'
$ws.Range("D46").Value = 'cut-in'
$ws.Range("E46").Value = 'm/s'
$ws.Range("F46").Value = 'moderate'
$ws.Range("G46").Value = 'acceptable'
$ws.Range("I46").Value = 'Cut-in wind speed. Source: WindTurbine power curves database.'
$ws.Range("J46").Value = 'triangular'
$ws.Range("K46").Value = 3.6
$ws.Range("L46").Value = 1.5
$ws.Range("M46").Value = 6.5
$ws.Range("N46").Value = 3.6
$ws.Range("O46").Value = 1.5
$ws.Range("P46").Value = 6.5
$ws.Range("Q46").Value = 3.6
$ws.Range("R46").Value = 1.5
$ws.Range("S46").Value = 6.5
$ws.Range("T46").Value = 3.6
$ws.Range("U46").Value = 1.5
$ws.Range("V46").Value = 6.5
$ws.Range("W46").Value = 3.6
$ws.Range("X46").Value = 1.5
$ws.Range("Y46").Value = 6.5
$ws.Range("Z46").Value = 3.6
$ws.Range("AA46").Value = 1.5
$ws.Range("AB46").Value = 6.5

# --- Row 47 ---
$ws.Range("A47").Value = 'Operation'
$ws.Range("B47").Value = 'all'
$ws.Range("C47").Value = 500
$ws.Range("C47").NumberFormat = "0"
$ws.Range("C47").HorizontalAlignment = -4108
$ws.Range("C47").VerticalAlignment = -4108
$ws.Range("D47").Value = 'cut-in'
$ws.Range("E47").Value = 'm/s'
$ws.Range("F47").Value = 'moderate'
$ws.Range("G47").Value = 'acceptable'
$ws.Range("I47").Value = 'Cut-in wind speed. Source: WindTurbine power curves database.'
$ws.Range("J47").Value = 'triangular'
$ws.Range("K47").Value = 3.6
$ws.Range("L47").Value = 1.5
$ws.Range("M47").Value = 6.5
$ws.Range("N47").Value = 3.6
$ws.Range("O47").Value = 1.5
$ws.Range("P47").Value = 6.5
$ws.Range("Q47").Value = 3.6
$ws.Range("R47").Value = 1.5
$ws.Range("S47").Value = 6.5
$ws.Range("T47").Value = 3.6
$ws.Range("U47").Value = 1.5
$ws.Range("V47").Value = 6.5
$ws.Range("W47").Value = 3.6
$ws.Range("X47").Value = 1.5
$ws.Range("Y47").Value = 6.5
$ws.Range("Z47").Value = 3.6
$ws.Range("AA47").Value = 1.5
$ws.Range("AB47").Value = 6.5

# --- Row 48 ---
$ws.Range("A48").Value = 'Operation'
$ws.Range("B48").Value = 'all'
$ws.Range("C48").Value = 1000
$ws.Range("C48").NumberFormat = "0"
$ws.Range("C48").HorizontalAlignment = -4108
$ws.Range("C48").VerticalAlignment = -4108
$ws.Range("D48").Value = 'cut-in'
$ws.Range("E48").Value = 'm/s'
$ws.Range("F48").Value = 'moderate'
$ws.Range("G48").Value = 'acceptable'
$ws.Range("I48").Value = 'Cut-in wind speed. Source: WindTurbine power curves database.'
$ws.Range("J48").Value = 'triangular'
$ws.Range("K48").Value = 3.3
$ws.Range("L48").Value = 2
$ws.Range("M48").Value = 5.5
$ws.Range("N48").Value = 3.3
$ws.Range("O48").Value = 2
$ws.Range("P48").Value = 5.5
$ws.Range("Q48").Value = 3.3
$ws.Range("R48").Value = 2
$ws.Range("S48").Value = 5.5
$ws.Range("T48").Value = 3.3
$ws.Range("U48").Value = 2
$ws.Range("V48").Value = 5.5
$ws.Range("W48").Value = 3.3
$ws.Range("X48").Value = 2
$ws.Range("Y48").Value = 5.5
$ws.Range("Z48").Value = 3.3
$ws.Range("AA48").Value = 2
$ws.Range("AB48").Value = 5.5

# --- Row 49 ---
$ws.Range("A49").Value = 'Operation'
$ws.Range("B49").Value = 'all'
$ws.Range("C49").Value = 3000
$ws.Range("C49").NumberFormat = "0"
$ws.Range("C49").HorizontalAlignment = -4108
$ws.Range("C49").VerticalAlignment = -4108
$ws.Range("D49").Value = 'cut-in'
$ws.Range("E49").Value = 'm/s'
$ws.Range("F49").Value = 'moderate'
$ws.Range("G49").Value = 'acceptable'
$ws.Range("I49").Value = 'Cut-in wind speed. Source: WindTurbine power curves database.'
$ws.Range("J49").Value = 'triangular'
$ws.Range("K49").Value = 3.1
$ws.Range("L49").Value = 2
$ws.Range("M49").Value = 5
$ws.Range("N49").Value = 3.1
$ws.Range("O49").Value = 2
$ws.Range("P49").Value = 5
$ws.Range("Q49").Value = 3.1
$ws.Range("R49").Value = 2
$ws.Range("S49").Value = 5
$ws.Range("T49").Value = 3.1
$ws.Range("U49").Value = 2
$ws.Range("V49").Value = 5
$ws.Range("W49").Value = 3.1
$ws.Range("X49").Value = 2
$ws.Range("Y49").Value = 5
$ws.Range("Z49").Value = 3.1
$ws.Range("AA49").Value = 2
$ws.Range("AB49").Value = 5

# --- Row 50 ---
$ws.Range("A50").Value = 'Operation'
$ws.Range("B50").Value = 'all'
$ws.Range("C50").Value = 8000
$ws.Range("C50").NumberFormat = "0"
$ws.Range("C50").HorizontalAlignment = -4108
$ws.Range("C50").VerticalAlignment = -4108
$ws.Range("D50").Value = 'cut-in'
$ws.Range("E50").Value = 'm/s'
$ws.Range("F50").Value = 'moderate'
$ws.Range("G50").Value = 'acceptable'
$ws.Range("I50").Value = 'Cut-in wind speed. Source: WindTurbine power curves database.'
$ws.Range("J50").Value = 'triangular'
$ws.Range("K50").Value = 3.1
$ws.Range("L50").Value = 2.5
$ws.Range("M50").Value = 4
$ws.Range("N50").Value = 3.1
$ws.Range("O50").Value = 2.5
$ws.Range("P50").Value = 4
$ws.Range("Q50").Value = 3.1
$ws.Range("R50").Value = 2.5
$ws.Range("S50").Value = 4
$ws.Range("T50").Value = 3.1
$ws.Range("U50").Value = 2.5
$ws.Range("V50").Value = 4
$ws.Range("W50").Value = 3.1
$ws.Range("X50").Value = 2.5
$ws.Range("Y50").Value = 4
$ws.Range("Z50").Value = 3.1
$ws.Range("AA50").Value = 2.5
$ws.Range("AB50").Value = 4

# --- Row 51 ---
$ws.Range("A51").Value = 'Operation'
$ws.Range("B51").Value = 'all'
$ws.Range("C51").Value = 100
$ws.Range("C51").NumberFormat = "0"
$ws.Range("C51").HorizontalAlignment = -4108
$ws.Range("C51").VerticalAlignment = -4108
$ws.Range("D51").Value = 'cut-out'
$ws.Range("E51").Value = 'm/s'
$ws.Range("F51").Value = 'moderate'
$ws.Range("G51").Value = 'acceptable'
$ws.Range("I51").Value = 'Cut-in wind speed. Source: WindTurbine power curves database.'
$ws.Range("J51").Value = 'triangular'
$ws.Range("L51").Value = 14
$ws.Range("M51").Value = 30
$ws.Range("N51").Value = 23.669871794871799
$ws.Range("N51").NumberFormat = "0.0"
$ws.Range("N51").HorizontalAlignment = -4108
$ws.Range("N51").VerticalAlignment = -4108
$ws.Range("O51").Value = 14
$ws.Range("P51").Value = 30
$ws.Range("Q51").Value = 23.669871794871799
$ws.Range("Q51").NumberFormat = "0.0"
$ws.Range("Q51").HorizontalAlignment = -4108
$ws.Range("Q51").VerticalAlignment = -4108
$ws.Range("R51").Value = 14
$ws.Range("S51").Value = 30
$ws.Range("T51").Value = 23.669871794871799
$ws.Range("T51").NumberFormat = "0.0"
$ws.Range("T51").HorizontalAlignment = -4108
$ws.Range("T51").VerticalAlignment = -4108
$ws.Range("U51").Value = 14
$ws.Range("V51").Value = 30
$ws.Range("W51").Value = 23.669871794871799
$ws.Range("W51").NumberFormat = "0.0"
$ws.Range("W51").HorizontalAlignment = -4108
$ws.Range("W51").VerticalAlignment = -4108
$ws.Range("X51").Value = 14
$ws.Range("Y51").Value = 30
$ws.Range("Z51").Value = 23.669871794871799
$ws.Range("Z51").NumberFormat = "0.0"
$ws.Range("Z51").HorizontalAlignment = -4108
$ws.Range("Z51").VerticalAlignment = -4108
$ws.Range("AA51").Value = 14
$ws.Range("AB51").Value = 30

# --- Row 52 ---
$ws.Range("A52").Value = 'Operation'
$ws.Range("B52").Value = 'all'
$ws.Range("C52").Value = 500
$ws.Range("C52").NumberFormat = "0"
$ws.Range("C52").HorizontalAlignment = -4108
$ws.Range("C52").VerticalAlignment = -4108
$ws.Range("D52").Value = 'cut-out'
$ws.Range("E52").Value = 'm/s'
$ws.Range("F52").Value = 'moderate'
$ws.Range("G52").Value = 'acceptable'
$ws.Range("I52").Value = 'Cut-in wind speed. Source: WindTurbine power curves database.'
$ws.Range("J52").Value = 'triangular'
$ws.Range("K52").Value = 23.669871794871799
$ws.Range("K52").NumberFormat = "0.0"
$ws.Range("K52").HorizontalAlignment = -4108
$ws.Range("K52").VerticalAlignment = -4108
$ws.Range("L52").Value = 14
$ws.Range("M52").Value = 30
$ws.Range("N52").Value = 23.669871794871799
$ws.Range("N52").NumberFormat = "0.0"
$ws.Range("N52").HorizontalAlignment = -4108
$ws.Range("N52").VerticalAlignment = -4108
$ws.Range("O52").Value = 14
$ws.Range("P52").Value = 30
$ws.Range("Q52").Value = 23.669871794871799
$ws.Range("Q52").NumberFormat = "0.0"
$ws.Range("Q52").HorizontalAlignment = -4108
$ws.Range("Q52").VerticalAlignment = -4108
$ws.Range("R52").Value = 14
$ws.Range("S52").Value = 30
$ws.Range("T52").Value = 23.669871794871799
$ws.Range("T52").NumberFormat = "0.0"
$ws.Range("T52").HorizontalAlignment = -4108
$ws.Range("T52").VerticalAlignment = -4108
$ws.Range("U52").Value = 14
$ws.Range("V52").Value = 30
$ws.Range("W52").Value = 23.669871794871799
$ws.Range("W52").NumberFormat = "0.0"
$ws.Range("W52").HorizontalAlignment = -4108
$ws.Range("W52").VerticalAlignment = -4108
$ws.Range("X52").Value = 14
$ws.Range("Y52").Value = 30
$ws.Range("Z52").Value = 23.669871794871799
$ws.Range("Z52").NumberFormat = "0.0"
$ws.Range("Z52").HorizontalAlignment = -4108
$ws.Range("Z52").VerticalAlignment = -4108
$ws.Range("AA52").Value = 14
$ws.Range("AB52").Value = 30

# --- Row 53 ---
$ws.Range("A53").Value = 'Operation'
$ws.Range("B53").Value = 'all'
$ws.Range("C53").Value = 1000
$ws.Range("C53").NumberFormat = "0"
$ws.Range("C53").HorizontalAlignment = -4108
$ws.Range("C53").VerticalAlignment = -4108
$ws.Range("D53").Value = 'cut-out'
$ws.Range("E53").Value = 'm/s'
$ws.Range("F53").Value = 'moderate'
$ws.Range("G53").Value = 'acceptable'
$ws.Range("I53").Value = 'Cut-in wind speed. Source: WindTurbine power curves database.'
$ws.Range("J53").Value = 'triangular'
$ws.Range("K53").Value = 23.675000000000001
$ws.Range("K53").NumberFormat = "0.0"
$ws.Range("K53").HorizontalAlignment = -4108
$ws.Range("K53").VerticalAlignment = -4108
$ws.Range("L53").Value = 15
$ws.Range("M53").Value = 25
$ws.Range("N53").Value = 23.675000000000001
$ws.Range("N53").NumberFormat = "0.0"
$ws.Range("N53").HorizontalAlignment = -4108
$ws.Range("N53").VerticalAlignment = -4108
$ws.Range("O53").Value = 15
$ws.Range("P53").Value = 25
$ws.Range("Q53").Value = 23.675000000000001
$ws.Range("Q53").NumberFormat = "0.0"
$ws.Range("Q53").HorizontalAlignment = -4108
$ws.Range("Q53").VerticalAlignment = -4108
$ws.Range("R53").Value = 15
$ws.Range("S53").Value = 25
$ws.Range("T53").Value = 23.675000000000001
$ws.Range("T53").NumberFormat = "0.0"
$ws.Range("T53").HorizontalAlignment = -4108
$ws.Range("T53").VerticalAlignment = -4108
$ws.Range("U53").Value = 15
$ws.Range("V53").Value = 25
$ws.Range("W53").Value = 23.675000000000001
$ws.Range("W53").NumberFormat = "0.0"
$ws.Range("W53").HorizontalAlignment = -4108
$ws.Range("W53").VerticalAlignment = -4108
$ws.Range("X53").Value = 15
$ws.Range("Y53").Value = 25
$ws.Range("Z53").Value = 23.675000000000001
$ws.Range("Z53").NumberFormat = "0.0"
$ws.Range("Z53").HorizontalAlignment = -4108
$ws.Range("Z53").VerticalAlignment = -4108
$ws.Range("AA53").Value = 15
$ws.Range("AB53").Value = 25

# --- Row 54 ---
$ws.Range("A54").Value = 'Operation'
$ws.Range("B54").Value = 'all'
$ws.Range("C54").Value = 3000
$ws.Range("C54").NumberFormat = "0"
$ws.Range("C54").HorizontalAlignment = -4108
$ws.Range("C54").VerticalAlignment = -4108
$ws.Range("D54").Value = 'cut-out'
$ws.Range("E54").Value = 'm/s'
$ws.Range("F54").Value = 'moderate'
$ws.Range("G54").Value = 'acceptable'
$ws.Range("I54").Value = 'Cut-in wind speed. Source: WindTurbine power curves database.'
$ws.Range("J54").Value = 'triangular'
$ws.Range("K54").Value = 23.79
$ws.Range("K54").NumberFormat = "0.0"
$ws.Range("K54").HorizontalAlignment = -4108
$ws.Range("K54").VerticalAlignment = -4108
$ws.Range("L54").Value = 20
$ws.Range("M54").Value = 25
$ws.Range("N54").Value = 23.79
$ws.Range("N54").NumberFormat = "0.0"
$ws.Range("N54").HorizontalAlignment = -4108
$ws.Range("N54").VerticalAlignment = -4108
$ws.Range("O54").Value = 20
$ws.Range("P54").Value = 25
$ws.Range("Q54").Value = 23.79
$ws.Range("Q54").NumberFormat = "0.0"
$ws.Range("Q54").HorizontalAlignment = -4108
$ws.Range("Q54").VerticalAlignment = -4108
$ws.Range("R54").Value = 20
$ws.Range("S54").Value = 25
$ws.Range("T54").Value = 23.79
$ws.Range("T54").NumberFormat = "0.0"
$ws.Range("T54").HorizontalAlignment = -4108
$ws.Range("T54").VerticalAlignment = -4108
$ws.Range("U54").Value = 20
$ws.Range("V54").Value = 25
$ws.Range("W54").Value = 23.79
$ws.Range("W54").NumberFormat = "0.0"
$ws.Range("W54").HorizontalAlignment = -4108
$ws.Range("W54").VerticalAlignment = -4108
$ws.Range("X54").Value = 20
$ws.Range("Y54").Value = 25
$ws.Range("Z54").Value = 23.79
$ws.Range("Z54").NumberFormat = "0.0"
$ws.Range("Z54").HorizontalAlignment = -4108
$ws.Range("Z54").VerticalAlignment = -4108
$ws.Range("AA54").Value = 20
$ws.Range("AB54").Value = 25

# --- Row 55 ---
$ws.Range("A55").Value = 'Operation'
$ws.Range("B55").Value = 'all'
$ws.Range("C55").Value = 8000
$ws.Range("C55").NumberFormat = "0"
$ws.Range("C55").HorizontalAlignment = -4108
$ws.Range("C55").VerticalAlignment = -4108
$ws.Range("D55").Value = 'cut-out'
$ws.Range("E55").Value = 'm/s'
$ws.Range("F55").Value = 'moderate'
$ws.Range("G55").Value = 'acceptable'
$ws.Range("I55").Value = 'Cut-in wind speed. Source: WindTurbine power curves database.'
$ws.Range("J55").Value = 'triangular'
$ws.Range("K55").Value = 24.875
$ws.Range("K55").NumberFormat = "0.0"
$ws.Range("K55").HorizontalAlignment = -4108
$ws.Range("K55").VerticalAlignment = -4108
$ws.Range("L55").Value = 24
$ws.Range("M55").Value = 25
$ws.Range("N55").Value = 24.875
$ws.Range("N55").NumberFormat = "0.0"
$ws.Range("N55").HorizontalAlignment = -4108
$ws.Range("N55").VerticalAlignment = -4108
$ws.Range("O55").Value = 24
$ws.Range("P55").Value = 25
$ws.Range("Q55").Value = 24.875
$ws.Range("Q55").NumberFormat = "0.0"
$ws.Range("Q55").HorizontalAlignment = -4108
$ws.Range("Q55").VerticalAlignment = -4108
$ws.Range("R55").Value = 24
$ws.Range("S55").Value = 25
$ws.Range("T55").Value = 24.875
$ws.Range("T55").NumberFormat = "0.0"
$ws.Range("T55").HorizontalAlignment = -4108
$ws.Range("T55").VerticalAlignment = -4108
$ws.Range("U55").Value = 24
$ws.Range("V55").Value = 25
$ws.Range("W55").Value = 24.875
$ws.Range("W55").NumberFormat = "0.0"
$ws.Range("W55").HorizontalAlignment = -4108
$ws.Range("W55").VerticalAlignment = -4108
$ws.Range("X55").Value = 24
$ws.Range("Y55").Value = 25
$ws.Range("Z55").Value = 24.875
$ws.Range("Z55").NumberFormat = "0.0"
$ws.Range("Z55").HorizontalAlignment = -4108
$ws.Range("Z55").VerticalAlignment = -4108
$ws.Range("AA55").Value = 24
$ws.Range("AB55").Value = 25

# Restore the selection to match the authored view (row 48, column D).
$ws.Range("D48").Select()
